# Update the three input cells on the active sheet ("Fin Buff Calc").
# D6:D9 (and E8:E9) are formulas that depend on these and recalc automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 340569.79
$ws.Range("D4").Value = 385.14
$ws.Range("D5").Value = 72315.09

# Reflect the saved selection state (merged range G7:I16) from the source file.
$ws.Range("G7:I16").Select()
